# Add new booking row (SNOW-777038) to the Bookings sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 7

$ws.Range("A$newRow").Value = "SNOW-777038"
# Leading apostrophe forces these to stay literal text instead of being
# auto-converted to a date/number by Excel (matches columns B/E/K on the
# existing rows, which are stored as text too).
$ws.Range("B$newRow").Value = "'2026-02-20"
$ws.Range("C$newRow").Value = "test"
$ws.Range("D$newRow").Value = "test@test.com"
$ws.Range("E$newRow").Value = "'1213454"
$ws.Range("F$newRow").Value = 5
$ws.Range("G$newRow").Value = "Family Ski Package"
$ws.Range("H$newRow").Value = 32000
$ws.Range("I$newRow").Value = 32000
$ws.Range("J$newRow").Value = "Confirmed"
$ws.Range("K$newRow").Value = "'2026-02-16"
$ws.Range("L$newRow").Value = "'"
